# Change the "CAREER PROGRESSION:" heading on slide 1 from the themed
# accent5 (lumMod 50%) fill to the explicit navy RGB color 002060, to
# match the rest of the sidebar headings/text in the same text box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$targetColor = 0x602000

$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "CAREER PROGRESSION:*") {
        $runCount = $para.Runs().Count
        for ($j = 1; $j -le $runCount; $j++) {
            $run = $para.Runs($j, 1)
            $run.Font.Color.RGB = $targetColor
        }
    }
}
